# "Updated Master data as per 16th May Refresh"
#
# This script replays, via Excel COM automation, the edits made to
# master-template_type.xlsx:
#   1. Renames three existing master-data codes (the "otp-*" templates used
#      by the IDA auth flow) to their new "ida-auth-otp-*-template" names,
#      in rows 5/6/7 (eng), 11/12/13 (ara) and 17/18/19 (fra).
#   2. Appends a new "consent" master-data row (eng/ara/fra) and three new
#      "auth-otp-*-template" master-data rows (eng/ara/fra each) at the
#      bottom of the table (rows 125-136).
#   3. Leaves the selection positioned just below the new last row, as
#      Excel does after editing the end of a used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the "otp-*" codes to "ida-auth-otp-*-template" -------------
# eng
$ws.Range("A7").Value = "ida-auth-otp-sms-template"
$ws.Range("A6").Value = "ida-auth-otp-email-subject-template"
$ws.Range("A5").Value = "ida-auth-otp-email-content-template"
# ara
$ws.Range("A13").Value = "ida-auth-otp-sms-template"
$ws.Range("A12").Value = "ida-auth-otp-email-subject-template"
$ws.Range("A11").Value = "ida-auth-otp-email-content-template"
# fra
$ws.Range("A19").Value = "ida-auth-otp-sms-template"
$ws.Range("A18").Value = "ida-auth-otp-email-subject-template"
$ws.Range("A17").Value = "ida-auth-otp-email-content-template"

# --- 2a. New "consent" rows (125-127) --------------------------------------
$ws.Range("A125").Value = "consent"
$ws.Range("B125").Value = "Consent"
$ws.Range("C125").Value = "eng"
$ws.Range("D125").Value = $true
$ws.Range("E125").Value = "superadmin"
$ws.Range("F125").Value = "now()"

$ws.Range("A126").Value = "consent"
$ws.Range("B126").Value = "موافقة"
$ws.Range("C126").Value = "ara"
$ws.Range("D126").Value = $true
$ws.Range("E126").Value = "superadmin"
$ws.Range("F126").Value = "now()"

$ws.Range("A127").Value = "consent"
$ws.Range("B127").Value = "Consentement"
$ws.Range("C127").Value = "fra"
$ws.Range("D127").Value = $true
$ws.Range("E127").Value = "superadmin"
$ws.Range("F127").Value = "now()"

# --- 2b. New "auth-otp-*-template" rows (128-136) --------------------------
# column A (codes) first ...
$ws.Range("A128").Value = "auth-otp-email-subject-template"
$ws.Range("A131").Value = "auth-otp-email-content-template"
$ws.Range("A134").Value = "auth-otp-sms-template"
$ws.Range("A129").Value = "auth-otp-email-subject-template"
$ws.Range("A130").Value = "auth-otp-email-subject-template"
$ws.Range("A132").Value = "auth-otp-email-content-template"
$ws.Range("A133").Value = "auth-otp-email-content-template"
$ws.Range("A135").Value = "auth-otp-sms-template"
$ws.Range("A136").Value = "auth-otp-sms-template"

# ... then lang_code / is_active / cr_by / cr_dtimes for each row ...
$ws.Range("C128").Value = "eng"
$ws.Range("D128").Value = $true
$ws.Range("E128").Value = "superadmin"
$ws.Range("F128").Value = "now()"

$ws.Range("C129").Value = "ara"
$ws.Range("D129").Value = $true
$ws.Range("E129").Value = "superadmin"
$ws.Range("F129").Value = "now()"

$ws.Range("C130").Value = "fra"
$ws.Range("D130").Value = $true
$ws.Range("E130").Value = "superadmin"
$ws.Range("F130").Value = "now()"

$ws.Range("C131").Value = "eng"
$ws.Range("D131").Value = $true
$ws.Range("E131").Value = "superadmin"
$ws.Range("F131").Value = "now()"

$ws.Range("C132").Value = "ara"
$ws.Range("D132").Value = $true
$ws.Range("E132").Value = "superadmin"
$ws.Range("F132").Value = "now()"

$ws.Range("C133").Value = "fra"
$ws.Range("D133").Value = $true
$ws.Range("E133").Value = "superadmin"
$ws.Range("F133").Value = "now()"

$ws.Range("C134").Value = "eng"
$ws.Range("D134").Value = $true
$ws.Range("E134").Value = "superadmin"
$ws.Range("F134").Value = "now()"

$ws.Range("C135").Value = "ara"
$ws.Range("D135").Value = $true
$ws.Range("E135").Value = "superadmin"
$ws.Range("F135").Value = "now()"

$ws.Range("C136").Value = "fra"
$ws.Range("D136").Value = $true
$ws.Range("E136").Value = "superadmin"
$ws.Range("F136").Value = "now()"

# ... finally column B (descriptions), in the order they were authored.
$ws.Range("B128").Value = "Auth OTP Email Subject Template"
$ws.Range("B131").Value = "Auth OTP Email Content Template"
$ws.Range("B134").Value = "Auth OTP SMS Template"
$ws.Range("B130").Value = "Modèle dobjet de-mail Auth OTP"
$ws.Range("B133").Value = "Auth OTP Email ContentTemplate"
$ws.Range("B136").Value = "Modèle SMS OTP Auth"
$ws.Range("B135").Value = "مصادقة قالب رسالة OTP"
$ws.Range("B132").Value = "مصادقة OTP قالب محتوى"
$ws.Range("B129").Value = "مصادقة OTP قالب موضوع"

# --- 3. Move the "next empty row" selection down to below the new data ----
$ws.Range("A137:XFD1048576").Select()
